$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$oldText = $cellA1.Value2
$newText = $oldText -replace [regex]::Escape("✅ 1000 Bs = 6.62 = 26165.56 pesos`n✅ 26165.56 pesos = 6.58 = 961.91 Bs"), "✅ 1000 Bs = 6.53 = 25905.65 pesos`n✅ 25905.65 pesos = 6.51 = 970.37 Bs"
$cellA1.Value2 = $newText

# --- Sheet "tasas": update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 153.094
$wsTasas.Range("O10").Value2 = 3966
$wsTasas.Range("N12").Value2 = 3977.8
$wsTasas.Range("O12").Value2 = 149
